# This script applies the "storeKeys" JSON command addition and the
# removal of the "text" target/category from the hidden '#system' sheet
# of the redis-showcase workbook.
#
# Summary of edits:
#   1. Remove the "text" entry (row 25) from the "target" list in column A,
#      shifting the rows below it (web, webalert, webcookie, ws, ws.async,
#      xml) up by one row. This is done with plain value assignment
#      (instead of Range.Delete) so that only column A is affected and the
#      other, unrelated columns sharing those row numbers are left intact.
#   2. Remove the entire "text" column (column Y), which only contained a
#      single value, shifting columns Z:AE left by one column to Y:AD.
#   3. Insert a new row for "storeKeys(json,jsonpath,var)" into the "json"
#      column (M) right before "storeValue(json,jsonpath,var)", shifting
#      storeValue/storeValues down by one row. Again done with plain value
#      assignment so only column M is touched.
#   4. Update the named ranges that are affected by the above shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. Remove the "text" entry from the target list (column A, row 25) ---
# Shift A26:A31 up to A25:A30 and clear the now-duplicated last row (A31).
for ($r = 25; $r -le 30; $r++) {
    $ws.Range("A$r").Value = $ws.Range("A$($r + 1)").Value()
}
$ws.Range("A31").Value = ""

# --- 2. Delete the now-orphaned "text" column (column Y) entirely ---
# Shifts Z:AE left to Y:AD.
$ws.Range("Y:Y").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# --- 3. Insert the new "storeKeys(json,jsonpath,var)" row into the json list ---
# Shift M16:M17 down to M17:M18 and set the freed M16 to the new entry.
$ws.Range("M18").Value = $ws.Range("M17").Value()
$ws.Range("M17").Value = $ws.Range("M16").Value()
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# --- 4. Fix up the defined (named) ranges impacted by the shifts above ---
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
